$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.911813333333333
$ws.Range("N2").Value = 5.73544
$ws.Range("O2").Value = 0.0369891942311654
$ws.Range("P2").Value = 0.03698919423116541
$ws.Range("Q2").Value = 1.045990673662222
$ws.Range("R2").Value = 9.41391606296
$ws.Range("S2").Value = 0.0369891942311654
$ws.Range("T2").Value = 0.03698919423116541

# Row 3 updates
$ws.Range("O3").Value = 0.005213386576832793
$ws.Range("P3").Value = 0.005213386576832794
$ws.Range("S3").Value = 0.005213386576832793
$ws.Range("T3").Value = 0.005213386576832794

# Row 4 updates
$ws.Range("M4").Value = 49.50445433333334
$ws.Range("N4").Value = 148.513363
$ws.Range("O4").Value = 0.9577974191920018
$ws.Range("P4").Value = 0.9577974191920018
$ws.Range("Q4").Value = 27.08486055336856
$ws.Range("R4").Value = 243.763744980317
$ws.Range("S4").Value = 0.9577974191920018
$ws.Range("T4").Value = 0.9577974191920018
